$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H12").Value = 12999.5
$ws_ALC.Range("I12").Value = 13000
$ws_ALC.Range("J12").Value = 12999
$ws_ALC.Range("K12").Value = 13000
$ws_ALC.Range("L12").Value = 12999
$ws_ALC.Range("M12").Value = -12830
$ws_ALC.Range("N12").Value = -13339
$ws_ALC.Range("H15").Value = 1375.569
$ws_ALC.Range("I15").Value = 1375.569
$ws_ALC.Range("K15").Value = 4126.707
$ws_ALC.Range("M15").Value = -3957.707
$ws_ALC.Range("H17").Value = 1670.7234
$ws_ALC.Range("J17").Value = 1670.7234
$ws_ALC.Range("L17").Value = 5012.1702
$ws_ALC.Range("N17").Value = -5348.1702
$ws_ALC.Range("H64").Value = 4181.25
$ws_ALC.Range("I64").Value = 3499.5
$ws_ALC.Range("J64").Value = 6226.5
$ws_ALC.Range("K64").Value = 3499.5
$ws_ALC.Range("L64").Value = 6226.5
$ws_ALC.Range("M64").Value = -3251.5
$ws_ALC.Range("N64").Value = -6722.5
$ws_ALC.Range("H67").Value = 4181.25
$ws_ALC.Range("I67").Value = 3499.5
$ws_ALC.Range("J67").Value = 6226.5
$ws_ALC.Range("K67").Value = 3499.5
$ws_ALC.Range("L67").Value = 6226.5
$ws_ALC.Range("M67").Value = -2641.5
$ws_ALC.Range("N67").Value = -7942.5
$ws_ALC.Range("H70").Value = 85031.336
$ws_ALC.Range("J70").Value = 251470.75
$ws_ALC.Range("L70").Value = 754412.25
$ws_ALC.Range("N70").Value = -754952.25
$ws_ALC.Range("H73").Value = 85031.336
$ws_ALC.Range("J73").Value = 251470.75
$ws_ALC.Range("L73").Value = 754412.25
$ws_ALC.Range("N73").Value = -756284.25
$ws_ALC.Range("H86").Value = 4051381.5
$ws_ALC.Range("J86").Value = 8773996
$ws_ALC.Range("L86").Value = 8773996
$ws_ALC.Range("N86").Value = -8776242
$ws_ALC.Range("H89").Value = 4051381.5
$ws_ALC.Range("J89").Value = 8773996
$ws_ALC.Range("L89").Value = 43869980
$ws_ALC.Range("N89").Value = -43881212
$ws_ALC.Range("H98").Value = 3018.1304
$ws_ALC.Range("I98").Value = 2313.625
$ws_ALC.Range("J98").Value = 4628.4287
$ws_ALC.Range("K98").Value = 2313.625
$ws_ALC.Range("L98").Value = 4628.4287
$ws_ALC.Range("M98").Value = -815.625
$ws_ALC.Range("N98").Value = -7624.4287
$ws_ALC.Range("H99").Value = 1357.625
$ws_ALC.Range("I99").Value = 215.25
$ws_ALC.Range("K99").Value = 645.75
$ws_ALC.Range("M99").Value = 852.25
$ws_ALC.Range("H100").Value = 9337.177
$ws_ALC.Range("I100").Value = 2127.5715
$ws_ALC.Range("K100").Value = 2127.5715
$ws_ALC.Range("M100").Value = -1586.5715
$ws_ALC.Range("H107").Value = 44301.914
$ws_ALC.Range("I107").Value = 48419.24
$ws_ALC.Range("J107").Value = 1070
$ws_ALC.Range("K107").Value = 48419.24
$ws_ALC.Range("L107").Value = 1070
$ws_ALC.Range("M107").Value = -46499.24
$ws_ALC.Range("N107").Value = -4910
$ws_ALC.Range("H112").Value = 3168.5625
$ws_ALC.Range("I112").Value = 350
$ws_ALC.Range("K112").Value = 1050
$ws_ALC.Range("M112").Value = 58
$ws_ALC.Range("H116").Value = 7108.4546
$ws_ALC.Range("I116").Value = 6666.3335
$ws_ALC.Range("J116").Value = 7639
$ws_ALC.Range("K116").Value = 6666.3335
$ws_ALC.Range("L116").Value = 7639
$ws_ALC.Range("M116").Value = -3224.3335
$ws_ALC.Range("N116").Value = -14523
$ws_ALC.Range("H122").Value = 3018.1304
$ws_ALC.Range("I122").Value = 2313.625
$ws_ALC.Range("J122").Value = 4628.4287
$ws_ALC.Range("K122").Value = 6940.875
$ws_ALC.Range("L122").Value = 13885.2861
$ws_ALC.Range("M122").Value = -4490.875
$ws_ALC.Range("N122").Value = -18785.2861
$ws_ALC.Range("H132").Value = 2571.0823
$ws_ALC.Range("I132").Value = 2183.8225
$ws_ALC.Range("K132").Value = 6551.467500000001
$ws_ALC.Range("M132").Value = -4021.467500000001
$ws_ALC.Range("H135").Value = 954069.0600000001
$ws_ALC.Range("I135").Value = 1429960.8
$ws_ALC.Range("K135").Value = 12869647.2
$ws_ALC.Range("M135").Value = -12867112.2
$ws_ALC.Range("H137").Value = 490972.62
$ws_ALC.Range("I137").Value = 314699.16
$ws_ALC.Range("K137").Value = 944097.48
$ws_ALC.Range("M137").Value = -941547.48
$ws_ALC.Range("H138").Value = 4116.9575
$ws_ALC.Range("J138").Value = 5019.6597
$ws_ALC.Range("L138").Value = 15058.9791
$ws_ALC.Range("N138").Value = -25338.9791
$ws_ALC.Range("H141").Value = 1769.8918
$ws_ALC.Range("I141").Value = 810.2222
$ws_ALC.Range("K141").Value = 2430.6666
$ws_ALC.Range("M141").Value = 2749.3334
$ws_ARM.Range("H45").Value = 3863.0833
$ws_ARM.Range("I45").Value = 2959.8333
$ws_ARM.Range("K45").Value = 2959.8333
$ws_ARM.Range("M45").Value = -2582.8333
$ws_ARM.Range("H61").Value = 1241.3148
$ws_ARM.Range("I61").Value = 700.62
$ws_ARM.Range("K61").Value = 700.62
$ws_ARM.Range("M61").Value = -488.62
$ws_ARM.Range("H74").Value = 3135
$ws_ARM.Range("I74").Value = 2986.625
$ws_ARM.Range("K74").Value = 2986.625
$ws_ARM.Range("M74").Value = -2112.625
$ws_ARM.Range("H77").Value = 3135
$ws_ARM.Range("I77").Value = 2986.625
$ws_ARM.Range("K77").Value = 14933.125
$ws_ARM.Range("M77").Value = -10565.125
$ws_ARM.Range("H122").Value = 3433.1538
$ws_ARM.Range("I122").Value = 1888.1923
$ws_ARM.Range("K122").Value = 5664.5769
$ws_ARM.Range("M122").Value = -3214.5769
$ws_ARM.Range("H136").Value = 1241.3148
$ws_ARM.Range("I136").Value = 700.62
$ws_ARM.Range("K136").Value = 2101.86
$ws_ARM.Range("M136").Value = 448.1399999999999
$ws_BSM.Range("H80").Value = 1588.5
$ws_BSM.Range("I80").Value = 548.25
$ws_BSM.Range("K80").Value = 548.25
$ws_BSM.Range("M80").Value = 449.75
$ws_BSM.Range("H83").Value = 1588.5
$ws_BSM.Range("I83").Value = 548.25
$ws_BSM.Range("K83").Value = 2741.25
$ws_BSM.Range("M83").Value = 2250.75
$ws_BSM.Range("H94").Value = 11021.4
$ws_BSM.Range("I94").Value = 2554.5
$ws_BSM.Range("J94").Value = 16666
$ws_BSM.Range("K94").Value = 2554.5
$ws_BSM.Range("L94").Value = 16666
$ws_BSM.Range("M94").Value = -2103.5
$ws_BSM.Range("N94").Value = -17568
$ws_CRP.Range("H31").Value = 2339502.2
$ws_CRP.Range("I31").Value = 10000012
$ws_CRP.Range("K31").Value = 10000012
$ws_CRP.Range("M31").Value = -9999717
$ws_CRP.Range("H34").Value = 2339502.2
$ws_CRP.Range("I34").Value = 10000012
$ws_CRP.Range("K34").Value = 10000012
$ws_CRP.Range("M34").Value = -9999810
$ws_CRP.Range("H58").Value = 200181.17
$ws_CRP.Range("I58").Value = 315131.47
$ws_CRP.Range("K58").Value = 315131.47
$ws_CRP.Range("M58").Value = -314928.47
$ws_CRP.Range("H132").Value = 2507.7693
$ws_CRP.Range("I132").Value = 1857.254
$ws_CRP.Range("K132").Value = 5571.762
$ws_CRP.Range("M132").Value = -3041.762
$ws_CRP.Range("H134").Value = 403582.62
$ws_CRP.Range("I134").Value = 252894.7
$ws_CRP.Range("J134").Value = 1006334.4
$ws_CRP.Range("K134").Value = 758684.1000000001
$ws_CRP.Range("L134").Value = 3019003.2
$ws_CRP.Range("M134").Value = -756149.1000000001
$ws_CRP.Range("N134").Value = -3024073.2
$ws_CRP.Range("H135").Value = 59614.92
$ws_CRP.Range("J135").Value = 59614.92
$ws_CRP.Range("L135").Value = 59614.92
$ws_CRP.Range("N135").Value = -69754.92
$ws_CRP.Range("H136").Value = 200181.17
$ws_CRP.Range("I136").Value = 315131.47
$ws_CRP.Range("K136").Value = 945394.4099999999
$ws_CRP.Range("M136").Value = -942844.4099999999
$ws_CUL.Range("H4").Value = 2117418.2
$ws_CUL.Range("I4").Value = 784546.4
$ws_CUL.Range("K4").Value = 2353639.2
$ws_CUL.Range("M4").Value = -2353527.2
$ws_CUL.Range("H92").Value = 527059.0600000001
$ws_CUL.Range("I92").Value = 1111683.6
$ws_CUL.Range("J92").Value = 896.9
$ws_CUL.Range("K92").Value = 3335050.8
$ws_CUL.Range("L92").Value = 2690.7
$ws_CUL.Range("M92").Value = -3333802.8
$ws_CUL.Range("N92").Value = -5186.7
$ws_CUL.Range("H136").Value = 7741.75
$ws_CUL.Range("J136").Value = 11993
$ws_CUL.Range("L136").Value = 35979
$ws_CUL.Range("N136").Value = -46179
$ws_GSM.Range("H3").Value = 6162000.5
$ws_GSM.Range("I3").Value = 6933334.5
$ws_GSM.Range("J3").Value = 5005000
$ws_GSM.Range("K3").Value = 6933334.5
$ws_GSM.Range("L3").Value = 5005000
$ws_GSM.Range("M3").Value = -6933218.5
$ws_GSM.Range("N3").Value = -5005232
$ws_GSM.Range("H102").Value = 2050.3684
$ws_GSM.Range("I102").Value = 1104.6786
$ws_GSM.Range("J102").Value = 4698.3
$ws_GSM.Range("K102").Value = 1104.6786
$ws_GSM.Range("L102").Value = 4698.3
$ws_GSM.Range("M102").Value = 517.3214
$ws_GSM.Range("N102").Value = -7942.3
$ws_GSM.Range("H128").Value = 73642.86
$ws_GSM.Range("J128").Value = 74250
$ws_GSM.Range("L128").Value = 74250
$ws_GSM.Range("N128").Value = -84210
$ws_LTW.Range("H46").Value = 3961.2778
$ws_LTW.Range("I46").Value = 3063.7273
$ws_LTW.Range("J46").Value = 5371.7144
$ws_LTW.Range("K46").Value = 3063.7273
$ws_LTW.Range("L46").Value = 5371.7144
$ws_LTW.Range("M46").Value = -2875.7273
$ws_LTW.Range("N46").Value = -5747.7144
$ws_LTW.Range("H61").Value = 4715.5625
$ws_LTW.Range("I61").Value = 2762.5
$ws_LTW.Range("K61").Value = 2762.5
$ws_LTW.Range("M61").Value = -2560.5
$ws_LTW.Range("H68").Value = 78735
$ws_LTW.Range("I68").Value = 4933
$ws_LTW.Range("K68").Value = 4933
$ws_LTW.Range("M68").Value = -4184
$ws_LTW.Range("H71").Value = 78735
$ws_LTW.Range("I71").Value = 4933
$ws_LTW.Range("K71").Value = 24665
$ws_LTW.Range("M71").Value = -20921
$ws_LTW.Range("H113").Value = 4715.5625
$ws_LTW.Range("I113").Value = 2762.5
$ws_LTW.Range("K113").Value = 2762.5
$ws_LTW.Range("M113").Value = -592.5
$ws_WVR.Range("H122").Value = 25003564
$ws_WVR.Range("I122").Value = 52633260
$ws_WVR.Range("K122").Value = 157899780
$ws_WVR.Range("M122").Value = -157897330
$ws_WVR.Range("H136").Value = 235814.98
$ws_WVR.Range("I136").Value = 306730.75
$ws_WVR.Range("K136").Value = 920192.25
$ws_WVR.Range("M136").Value = -917642.25
